$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Backend Software Engineer", "https://www.dice.com/job-detail/2f164b47-a7b0-4780-afb0-9fbee67daf87?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang", "Remote or Evansville, Illinois", "Contract", "Market rate", "Envision Technology Solutions"),
    @("UX PM(Need candidates local to Chicago iL)", "https://www.dice.com/job-detail/017b6b39-b14f-4aaa-a8ea-56002f56464a?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=3&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang", "Hybrid in Chicago, Illinois", "Contract", "$65 - $70", "DTEL Engineering & Consultants Inc"),
    @("AWS Full-Stack Backend Developer", "https://www.dice.com/job-detail/ab21904b-4096-4b9c-83c6-b1adb83a7809?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=4&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang", "Albany, New York", "Full-time, Third Party, Contract", "$12 - $16", "Xoriant Corporation")
)

$startRow = 192
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
